# The workbook tracks daily price observations for "Coliflor" at the
# Macroferia Regional de Talca market. A new weekly observation is being
# inserted as row 62 (pushing the existing rows 62-144 down to 63-145).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 62, shifting rows 62:144 down to 63:145.
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new observation.
$ws.Cells.Item(62, 1).Value()  = 5
$ws.Cells.Item(62, 2).Value()  = "Macroferia Regional de Talca"
$ws.Cells.Item(62, 3).Value()  = "Maule"
$ws.Cells.Item(62, 4).Value()  = 44467
$ws.Cells.Item(62, 5).Value()  = 7
$ws.Cells.Item(62, 6).Value()  = 100112008
$ws.Cells.Item(62, 7).Value()  = "Coliflor"
$ws.Cells.Item(62, 8).Value()  = "Sin especificar"
$ws.Cells.Item(62, 9).Value()  = "Primera"
$ws.Cells.Item(62, 10).Value() = 3000
$ws.Cells.Item(62, 11).Value() = 600
$ws.Cells.Item(62, 12).Value() = 600
$ws.Cells.Item(62, 13).Value() = 600
$ws.Cells.Item(62, 14).Value() = "`$/unidad"
$ws.Cells.Item(62, 15).Value() = "Región del Maule"
$ws.Cells.Item(62, 16).Value() = 600
$ws.Cells.Item(62, 17).Value() = 1
$ws.Cells.Item(62, 18).Value() = "Hortaliza"
